# Auto-generated Excel COM-interop script
# Updates market-price-derived columns (H-N) across 8 sheets to reflect
# refreshed values from the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 735
$ws.Range("J2").Value = 735
$ws.Range("L2").Value = 735
$ws.Range("N2").Value = -961
$ws.Range("H32").Value = 15872
$ws.Range("I32").Value = 14244.5
$ws.Range("J32").Value = 17499.5
$ws.Range("K32").Value = 14244.5
$ws.Range("L32").Value = 17499.5
$ws.Range("M32").Value = -13918.5
$ws.Range("N32").Value = -18151.5
$ws.Range("H61").Value = 739.6
$ws.Range("I61").Value = 739.6
$ws.Range("K61").Value = 2218.8
$ws.Range("M61").Value = -2046.8
$ws.Range("H76").Value = 5165.8335
$ws.Range("J76").Value = 5498.75
$ws.Range("L76").Value = 5498.75
$ws.Range("N76").Value = -6128.75
$ws.Range("H79").Value = 5165.8335
$ws.Range("J79").Value = 5498.75
$ws.Range("L79").Value = 5498.75
$ws.Range("N79").Value = -7682.75
$ws.Range("H80").Value = 509.6875
$ws.Range("I80").Value = 202.14285
$ws.Range("J80").Value = 748.8889
$ws.Range("K80").Value = 606.4285500000001
$ws.Range("L80").Value = 2246.6667
$ws.Range("M80").Value = 391.5714499999999
$ws.Range("N80").Value = -4242.6667
$ws.Range("H82").Value = 11262.454
$ws.Range("I82").Value = 11262.454
$ws.Range("K82").Value = 33787.362
$ws.Range("M82").Value = -33381.362
$ws.Range("H83").Value = 509.6875
$ws.Range("I83").Value = 202.14285
$ws.Range("J83").Value = 748.8889
$ws.Range("K83").Value = 1819.28565
$ws.Range("L83").Value = 6740.0001
$ws.Range("M83").Value = 3172.71435
$ws.Range("N83").Value = -16724.0001
$ws.Range("H85").Value = 11262.454
$ws.Range("I85").Value = 11262.454
$ws.Range("K85").Value = 33787.362
$ws.Range("M85").Value = -32383.362
$ws.Range("H88").Value = 1534.95
$ws.Range("J88").Value = 1800.3
$ws.Range("L88").Value = 1800.3
$ws.Range("N88").Value = -2612.3
$ws.Range("H91").Value = 1534.95
$ws.Range("J91").Value = 1800.3
$ws.Range("L91").Value = 1800.3
$ws.Range("N91").Value = -4608.3
$ws.Range("H92").Value = 1415.7142
$ws.Range("I92").Value = 1222
$ws.Range("J92").Value = 1900
$ws.Range("K92").Value = 1222
$ws.Range("L92").Value = 1900
$ws.Range("M92").Value = 26
$ws.Range("N92").Value = -4396
$ws.Range("H101").Value = 12700.375
$ws.Range("I101").Value = 1156.091
$ws.Range("J101").Value = 38097.8
$ws.Range("K101").Value = 3468.273
$ws.Range("L101").Value = 114293.4
$ws.Range("M101").Value = -1846.273
$ws.Range("N101").Value = -117537.4
$ws.Range("H125").Value = 901.86664
$ws.Range("I125").Value = 784.1667
$ws.Range("J125").Value = 980.3333
$ws.Range("K125").Value = 7057.5003
$ws.Range("L125").Value = 8822.9997
$ws.Range("M125").Value = -4597.5003
$ws.Range("N125").Value = -13742.9997
$ws.Range("H127").Value = 1749.5
$ws.Range("I127").Value = 1714
$ws.Range("K127").Value = 5142
$ws.Range("M127").Value = -182
$ws.Range("H132").Value = 1813.2069
$ws.Range("I132").Value = 1566.875
$ws.Range("K132").Value = 4700.625
$ws.Range("M132").Value = -2170.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1999.875
$ws.Range("I74").Value = 1346.4286
$ws.Range("K74").Value = 1346.4286
$ws.Range("M74").Value = -472.4286
$ws.Range("H77").Value = 1999.875
$ws.Range("I77").Value = 1346.4286
$ws.Range("K77").Value = 6732.143
$ws.Range("M77").Value = -2364.143
$ws.Range("H132").Value = 1065.1515
$ws.Range("I132").Value = 914.09375
$ws.Range("J132").Value = 5899
$ws.Range("K132").Value = 2742.28125
$ws.Range("L132").Value = 17697
$ws.Range("M132").Value = -212.28125
$ws.Range("N132").Value = -22757

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8520.5
$ws.Range("I20").Value = 9871.777
$ws.Range("J20").Value = 4466.6665
$ws.Range("K20").Value = 9871.777
$ws.Range("L20").Value = 4466.6665
$ws.Range("M20").Value = -9624.777
$ws.Range("N20").Value = -4960.6665
$ws.Range("H39").Value = 44500
$ws.Range("J39").Value = 44500
$ws.Range("L39").Value = 44500
$ws.Range("N39").Value = -45278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3320
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726
$ws.Range("H29").Value = 15633
$ws.Range("J29").Value = 15633
$ws.Range("L29").Value = 15633
$ws.Range("N29").Value = -16219
$ws.Range("H105").Value = 1703.8889
$ws.Range("I105").Value = 1639.25
$ws.Range("J105").Value = 2221
$ws.Range("K105").Value = 1639.25
$ws.Range("L105").Value = 2221
$ws.Range("M105").Value = 107.75
$ws.Range("N105").Value = -5715
$ws.Range("H134").Value = 5506
$ws.Range("I134").Value = 4134
$ws.Range("K134").Value = 12402
$ws.Range("M134").Value = -9867

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1053387.4
$ws.Range("I4").Value = 1111353.5
$ws.Range("J4").Value = 9999
$ws.Range("K4").Value = 3334060.5
$ws.Range("L4").Value = 29997
$ws.Range("M4").Value = -3333948.5
$ws.Range("N4").Value = -30221
$ws.Range("H131").Value = 1574.0182
$ws.Range("I131").Value = 1115.2
$ws.Range("K131").Value = 3345.6
$ws.Range("M131").Value = 1694.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1148.6
$ws.Range("I80").Value = 998
$ws.Range("K80").Value = 998
$ws.Range("M80").Value = 0
$ws.Range("H83").Value = 1148.6
$ws.Range("I83").Value = 998
$ws.Range("K83").Value = 4990
$ws.Range("M83").Value = 2
$ws.Range("H113").Value = 3499.5
$ws.Range("J113").Value = 3499.5
$ws.Range("L113").Value = 3499.5
$ws.Range("N113").Value = -7839.5
$ws.Range("H122").Value = 3839.5
$ws.Range("I122").Value = 5982.8
$ws.Range("J122").Value = 2648.7778
$ws.Range("K122").Value = 17948.4
$ws.Range("L122").Value = 7946.3334
$ws.Range("M122").Value = -15498.4
$ws.Range("N122").Value = -12846.3334
$ws.Range("H132").Value = 2749
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 26202.79
$ws.Range("J136").Value = 26202.79
$ws.Range("L136").Value = 78608.37
$ws.Range("N136").Value = -83708.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3052.1365
$ws.Range("I46").Value = 1014
$ws.Range("J46").Value = 4003.2666
$ws.Range("K46").Value = 1014
$ws.Range("L46").Value = 4003.2666
$ws.Range("M46").Value = -826
$ws.Range("N46").Value = -4379.2666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4413.778
$ws.Range("I62").Value = 4873.25
$ws.Range("K62").Value = 4873.25
$ws.Range("M62").Value = -4249.25
$ws.Range("H65").Value = 4413.778
$ws.Range("I65").Value = 4873.25
$ws.Range("K65").Value = 24366.25
$ws.Range("M65").Value = -21246.25
$ws.Range("H122").Value = 615.9
$ws.Range("I122").Value = 569.875
$ws.Range("K122").Value = 1709.625
$ws.Range("M122").Value = 740.375
$ws.Range("H132").Value = 2613.5625
$ws.Range("I132").Value = 1987.5714
$ws.Range("K132").Value = 5962.7142
$ws.Range("M132").Value = -3432.7142
$ws.Range("H136").Value = 4004
$ws.Range("I136").Value = 4004
$ws.Range("K136").Value = 12012
$ws.Range("M136").Value = -9462
